$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "MNC Experiment" -> sheet1.xml
$ws2 = $wb.Worksheets.Item(2)   # "CAH Experiment" -> sheet2.xml

# ---------------------------------------------------------------------------
# Sheet1 ("MNC Experiment") - header row additions (split the combined label
# into separate factorybg/messaging/navprimary/knowledge columns) and
# highlight the STBO row (row 14).
# ---------------------------------------------------------------------------

# Row 7 header cells
$ws1.Range("C7").Interior.Color = 0x00FFFF
$ws1.Range("D7").Value = "response"
$ws1.Range("D7").Interior.Color = 0x00FFFF

$ws1.Range("E7").Value = "factorybg"
$ws1.Range("F7").Value = "messaging"
$ws1.Range("G7").Value = "navprimary"
$ws1.Range("H7").Value = "knowledge"
$ws1.Range("E7:H7").Interior.Color = 0xF0B000

# Row 14 (STBO result row) highlighting
$ws1.Range("C14").Font.Color = 0x0000FF
$ws1.Range("E14:H14").Interior.Color = 0x50D092

# New column widths for the newly-used columns E:H
$ws1.Columns.Item(5).ColumnWidth = 11.333333333333332
$ws1.Columns.Item(6).ColumnWidth = 11.333333333333332
$ws1.Columns.Item(7).ColumnWidth = 10.833333333333332
$ws1.Columns.Item(8).ColumnWidth = 11.666666666666666

# ---------------------------------------------------------------------------
# Sheet2 ("CAH Experiment") - re-style header row, drop obsolete BCBO row,
# rename/realign the baseline block.
# ---------------------------------------------------------------------------

$ws2.Range("G7").Interior.Color = 0x00FFFF
$ws2.Range("H7:K7").Interior.Color = 0xF0B000

# Rename the (lowercase) baseline label to "Baseline" before the row shift.
$ws2.Range("G33").Value = "Baseline"

# Remove the now-duplicated BCBO data row (row 30); everything below row 24
# shifts up by one, moving the BCBO label from row 25 to row 24 and the
# trailing baseline block from rows 33/34 to rows 32/33.
$ws2.Rows.Item(30).Delete()

# New column width for column G, and widen column J
$ws2.Columns.Item(7).ColumnWidth = 12.5
$ws2.Columns.Item(10).ColumnWidth = 14.333333333333332

# ---------------------------------------------------------------------------
# Selections / active sheet to match the saved view state.
# ---------------------------------------------------------------------------

$ws1.Range("M14").Select()
$ws2.Activate()
$ws2.Range("F10").Select()
